$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the values in columns B through H, rows 2-13 to the nearest integer,
# so that the Ontpl_/Pot_ reference files are written as integer data.
for ($r = 2; $r -le 13; $r++) {
    for ($c = 2; $c -le 8; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $current = [double]$cell.Value()
        $cell.Value = [math]::Round($current, 0)
    }
}
